$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4").Value = "5303710149825215"
$ws.Range("O2").Value = "4513070436920974"
$ws.Range("O5").Value = "5303710095505365"
$ws.Range("Q2").Value = "1000"
$ws.Range("Q3").Value = "1200"
